$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 58 (shifts existing rows 58-66 down to 60-68)
$ws.Rows.Item(58).Insert()
$ws.Rows.Item(58).Insert()

# Populate the first new row (58) with the latest weekly price data
$ws.Cells.Item(58, 1).Value = 1
$ws.Cells.Item(58, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(58, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(58, 4).Value = 44449
$ws.Cells.Item(58, 5).Value = 15
$ws.Cells.Item(58, 6).Value = 100112042
$ws.Cells.Item(58, 7).Value = "Locoto"
$ws.Cells.Item(58, 8).Value = "Sin especificar"
$ws.Cells.Item(58, 9).Value = "Primera"
$ws.Cells.Item(58, 10).Value = 120
$ws.Cells.Item(58, 11).Value = 24000
$ws.Cells.Item(58, 12).Value = 25000
$ws.Cells.Item(58, 13).Value = 24500
$ws.Cells.Item(58, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(58, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(58, 16).Value = 1225
$ws.Cells.Item(58, 17).Value = 20
$ws.Cells.Item(58, 18).Value = "Hortaliza"

# Populate the second new row (59) with the latest weekly price data
$ws.Cells.Item(59, 1).Value = 1
$ws.Cells.Item(59, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(59, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(59, 4).Value = 44449
$ws.Cells.Item(59, 5).Value = 15
$ws.Cells.Item(59, 6).Value = 100112042
$ws.Cells.Item(59, 7).Value = "Locoto"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Segunda"
$ws.Cells.Item(59, 10).Value = 160
$ws.Cells.Item(59, 11).Value = 22000
$ws.Cells.Item(59, 12).Value = 23000
$ws.Cells.Item(59, 13).Value = 22500
$ws.Cells.Item(59, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(59, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(59, 16).Value = 1125
$ws.Cells.Item(59, 17).Value = 20
$ws.Cells.Item(59, 18).Value = "Hortaliza"
